$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.552.93"
$ws.Range("E2").Value = "  +5.37%  "
$ws.Range("D3").Value = "2.253.56"
$ws.Range("E3").Value = "  +4.50%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'234.24"
$ws.Range("E5").Value = "  +2.90%  "
$ws.Range("D6").Value = "'0.637"
$ws.Range("E6").Value = "  +2.18%  "
$ws.Range("D7").Value = "'64.86"
$ws.Range("E7").Value = "  +1.33%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "'0.410"
$ws.Range("E9").Value = "  +3.39%  "
$ws.Range("D10").Value = "'59.62"
$ws.Range("E10").Value = "  +3.10%  "
$ws.Range("D11").Value = "'0.0904"
$ws.Range("E11").Value = "  +5.48%  "
$ws.Range("D12").Value = "'0.105"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").Value = "2.593.09"
$ws.Range("E13").Value = "  +4.64%  "
$ws.Range("D14").Value = "'16.19"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").Value = "'22.44"
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("E16").Value = "  +2.85%  "
$ws.Range("D17").Value = "'5.68"
$ws.Range("E17").Value = "  +2.53%  "
$ws.Range("D18").Value = "2.259.48"
$ws.Range("E18").Value = "  +3.87%  "
$ws.Range("D19").Value = "41.447.81"
$ws.Range("E19").Value = "  +5.29%  "
$ws.Range("D20").Value = "'74.12"
$ws.Range("E20").Value = "  +3.24%  "
$ws.Range("D21").Value = "0.0₃0918"
$ws.Range("E21").Value = "  +8.06%  "
$ws.Range("D22").Value = "'6.20"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("D23").Value = "'252.81"
$ws.Range("E23").Value = "  +9.51%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  +3.12%  "
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("D27").Value = "'9.83"
$ws.Range("E27").Value = "  +3.63%  "
$ws.Range("D28").Value = "'173.27"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("D29").Value = "'0.144"
$ws.Range("E29").Value = "  +3.07%  "
$ws.Range("E30").Value = "  +3.22%  "
$ws.Range("D31").Value = "'1.45"
$ws.Range("E31").Value = "  +2.59%  "
$ws.Range("D32").Value = "'2.81"
$ws.Range("E32").Value = "  +5.69%  "
$ws.Range("E33").Value = "  +2.55%  "
$ws.Range("D34").Value = "'4.74"
$ws.Range("E34").Value = "  +3.13%  "
$ws.Range("D35").Value = "'4.95"
$ws.Range("E35").Value = "  +4.48%  "
$ws.Range("D36").Value = "'7.20"
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("D37").Value = "'0.0636"
$ws.Range("E37").Value = "  +3.16%  "
$ws.Range("D38").Value = "'3.89"
$ws.Range("E38").Value = "  +8.91%  "
$ws.Range("D39").Value = "'2.46"
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("B41").Value = "TerraClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D41").Value = "'0.000236"
$ws.Range("E41").Value = "  +48.89%  "
$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D42").Value = "'4.87"
$ws.Range("E42").Value = "  +13.07%  "
$ws.Range("E43").Value = "  +3.47%  "
$ws.Range("D44").Value = "'8.81"
$ws.Range("E44").Value = "  +14.70%  "
$ws.Range("D45").Value = "'18.25"
$ws.Range("E45").Value = "  +3.66%  "
$ws.Range("D46").Value = "'102.33"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").Value = "'1.24"
$ws.Range("E47").Value = "  +4.28%  "
$ws.Range("D48").Value = "1.515.81"
$ws.Range("E48").Value = "  -1.56%  "
$ws.Range("D49").Value = "'0.0942"
$ws.Range("E49").Value = "  +1.18%  "
$ws.Range("D50").Value = "'1.11"
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("E51").Value = "  -0.38%  "
